$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.972.91'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.406.52'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.76'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.33'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  -4.45%  '
$ws.Range('D9').Value = '2.402.32'
$ws.Range('E9').Value = '  -3.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.104'
$ws.Range('E10').Value = '  -2.58%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.37'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.341'
$ws.Range('E13').Value = '  -4.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.41'
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('D15').Value = '2.848.09'
$ws.Range('E15').Value = '  -3.49%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '61.541.23'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000163'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '2.408.36'
$ws.Range('E18').Value = '  -3.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.62'
$ws.Range('E19').Value = '  -4.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.75'
$ws.Range('E20').Value = '  -3.86%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.84'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  -3.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.23'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').Value = '2.532.62'
$ws.Range('E27').Value = '  -3.64%  '
$ws.Range('D28').Value = '0.0₃0939'
$ws.Range('E28').Value = '  -8.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.67'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.43'
$ws.Range('E30').Value = '  -4.33%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.98'
$ws.Range('E31').Value = '  -4.34%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '509.93'
$ws.Range('E32').Value = '  -5.92%  '
$ws.Range('E33').Value = '  -2.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.49'
$ws.Range('E37').Value = '  -6.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.61'
$ws.Range('E38').Value = '  -5.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.373'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.00'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.50'
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.29'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.16'
$ws.Range('E45').Value = '  -6.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '140.33'
$ws.Range('E46').Value = '  -6.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.51'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.43'
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.578'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0919'
$ws.Range('E51').Value = '  -2.57%  '
